# "Kode pos" (column G) was holding the wrong data (phone-number-looking
# strings), while the real 5-digit postal codes were sitting in the
# "Nomor rumah" column (H). Fix column G with the correct postal codes and
# drop the now-redundant "Nomor rumah" column entirely (header + data),
# leaving "Status" where it already lives, in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the real postal codes (previously mislabeled under "Nomor rumah" /
# column H) into "Kode pos" (column G).
$ws.Range("G2").Value = 12345
$ws.Range("G3").Value = 12345
$ws.Range("G4").Value = 90909

# Remove the obsolete "Nomor rumah" column's header and data (column H).
# "Status" (column I) is left untouched in place.
$ws.Range("H1:H4").ClearContents()

# Let the now-stale column widths recompute for the columns whose content
# changed.
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()

Write-Output "G2:G4 = $($ws.Range('G2').Value2), $($ws.Range('G3').Value2), $($ws.Range('G4').Value2)"
Write-Output "H1 (should be empty) = [$($ws.Range('H1').Value2)]"
Write-Output "I1 = $($ws.Range('I1').Value2); I2 = $($ws.Range('I2').Value2)"
